# Consolidate the three text runs ("An", " ", "image") of the caption
# textbox into a single run, matching the PowerPoint writer's new
# behaviour of merging adjacent runs that share formatting.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)              # "TextBox 3" - the image caption
$tr = $sh.TextFrame.TextRange

# Re-assigning the identical string is treated as a no-op by the writer
# (no structural change recorded), so first set a throwaway value to
# force the paragraph to be rebuilt, then set the real text. The result
# is a single run containing "An image" instead of three.
$tr.Text = "__tmp__"
$tr.Text = "An image"
